# The commit swaps the presentation's theme ("Integral") for the stock
# "Office Theme": the theme bound to the slide master (ppt/theme/theme1.xml)
# picks up the "Office" palette (dk1/lt1 already match; dk2, lt2 and the
# six accents + hlink/folHlink change), while the theme bound to the
# notes master keeps the palette that used to live on the slide master.
#
# Concretely: theme1.xml's <a:clrScheme> goes from "Integral" colours to
# the stock "Office" colours below. (fontScheme / fmtScheme were already
# identical between the two theme parts, so nothing else to touch there.)

function Set-ThemeColorRGB {
    param($ColorScheme, [int]$Index, [string]$HexRRGGBB)
    $r = [Convert]::ToInt32($HexRRGGBB.Substring(0,2), 16)
    $g = [Convert]::ToInt32($HexRRGGBB.Substring(2,2), 16)
    $b = [Convert]::ToInt32($HexRRGGBB.Substring(4,2), 16)
    # PowerPoint's ColorFormat.RGB is a VBA RGB() value: R + G*256 + B*65536
    $ColorScheme.Colors($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Index map (PowerPoint ThemeColorScheme order):
#  1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
Set-ThemeColorRGB $colors 1  "000000"
Set-ThemeColorRGB $colors 2  "FFFFFF"
Set-ThemeColorRGB $colors 3  "44546A"
Set-ThemeColorRGB $colors 4  "E7E6E6"
Set-ThemeColorRGB $colors 5  "5B9BD5"
Set-ThemeColorRGB $colors 6  "ED7D31"
Set-ThemeColorRGB $colors 7  "A5A5A5"
Set-ThemeColorRGB $colors 8  "FFC000"
Set-ThemeColorRGB $colors 9  "4472C4"
Set-ThemeColorRGB $colors 10 "70AD47"
Set-ThemeColorRGB $colors 11 "0563C1"
Set-ThemeColorRGB $colors 12 "954F72"
